$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The quarterly rows for each year were recorded in the wrong B/C order
# (e.g. row 3 held "...B" data and row 4 held "...C" data). Swap the
# A:E content of each such pair back into the correct row positions.
$pairs = @(
    @(3, 4),
    @(7, 8),
    @(11, 12),
    @(15, 16),
    @(19, 20),
    @(23, 24),
    @(27, 28)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("A" + $r1 + ":E" + $r1)
    $range2 = $ws.Range("A" + $r2 + ":E" + $r2)

    $v1 = $range1.Value2
    $v2 = $range2.Value2

    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

# Columns F ("大型拖拉机产销率") and G ("大型拖拉机销售量") were a
# duplicate of earlier columns and are removed entirely.
$ws.Range("F1:G29").Delete()
